# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# on the active sheet to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new Price (column D) / Volume (column E) text. $null means "leave
# that column untouched" (only one of the two columns changed for that row).
$updates = @(
    @{ Row = 2;  D = "66.514.31"; E = "  -3.99%  " }
    @{ Row = 3;  D = "3.569.21";  E = "  -4.50%  " }
    @{ Row = 4;  D = $null;       E = "  -0.01%  " }
    @{ Row = 5;  D = "588.11";    E = "  -4.06%  " }
    @{ Row = 6;  D = "185.19";    E = "  -1.74%  " }
    @{ Row = 7;  D = "3.565.76";  E = "  -4.46%  " }
    @{ Row = 8;  D = "0.615";     E = "  -4.04%  " }
    @{ Row = 9;  D = $null;       E = "  +0.04%  " }
    @{ Row = 10; D = "0.672";     E = "  -7.15%  " }
    @{ Row = 11; D = "0.147";     E = "  -9.65%  " }
    @{ Row = 12; D = "53.26";     E = "  -7.46%  " }
    @{ Row = 13; D = $null;       E = "  -11.30%  " }
    @{ Row = 14; D = "9.84";      E = $null }
    @{ Row = 15; D = "4.135.52";  E = $null }
    @{ Row = 16; D = "3.565.00";  E = "  -4.60%  " }
    @{ Row = 17; D = $null;       E = "  -1.05%  " }
    @{ Row = 18; D = "18.37";     E = "  -5.34%  " }
    @{ Row = 19; D = "12.27";     E = $null }
    @{ Row = 20; D = "66.376.68"; E = "  -3.88%  " }
    @{ Row = 21; D = $null;       E = "  -7.36%  " }
    @{ Row = 22; D = "396.66";    E = "  -4.46%  " }
    @{ Row = 23; D = "4.37";      E = "  -5.77%  " }
    @{ Row = 24; D = "86.13";     E = "  -3.78%  " }
    @{ Row = 25; D = "11.43";     E = "  +3.97%  " }
    @{ Row = 26; D = "2.91";      E = "  -5.22%  " }
    @{ Row = 27; D = "12.50";     E = "  -3.35%  " }
    @{ Row = 28; D = $null;       E = "  -0.39%  " }
    @{ Row = 29; D = "3.55";      E = "  -6.63%  " }
    @{ Row = 30; D = "8.99";      E = "  -7.43%  " }
    @{ Row = 31; D = "31.18";     E = "  -6.67%  " }
    @{ Row = 32; D = "7.10";      E = "  -3.57%  " }
    @{ Row = 33; D = "12.20";     E = "  -4.73%  " }
    @{ Row = 34; D = "620.66";    E = "  -0.13%  " }
    @{ Row = 35; D = $null;       E = "  -7.63%  " }
    @{ Row = 36; D = "63.30";     E = "  -3.60%  " }
    @{ Row = 37; D = "41.49";     E = "  -7.88%  " }
    @{ Row = 38; D = $null;       E = "  -0.01%  " }
    @{ Row = 39; D = "0.404";     E = "  -2.37%  " }
    @{ Row = 40; D = $null;       E = "  -9.87%  " }
    @{ Row = 41; D = "0.132";     E = "  -6.68%  " }
    @{ Row = 42; D = "0.998";     E = "  -0.13%  " }
    @{ Row = 43; D = "3.004.15";  E = "  +5.97%  " }
    @{ Row = 44; D = $null;       E = "  -8.11%  " }
    @{ Row = 45; D = "2.54";      E = "  -4.38%  " }
    @{ Row = 46; D = $null;       E = "  -7.80%  " }
    @{ Row = 47; D = $null;       E = "  -7.67%  " }
    @{ Row = 48; D = "3.08";      E = "  -1.43%  " }
    @{ Row = 49; D = "8.59";      E = "  -7.12%  " }
    @{ Row = 50; D = "137.81";    E = "  -4.18%  " }
    @{ Row = 51; D = "2.75";      E = "  -1.43%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Prefix with an apostrophe so values that look numeric (e.g.
        # "588.11") are stored as plain text, matching the original
        # inline-string cells, instead of being auto-converted to numbers.
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
